# Adding field "type" to Contract
#
# Adds a new column O ("tipo") to the contracts sheet, populated for the
# first four data rows, widens column D, and moves the selection to the
# new last cell (O5) the way Excel would after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + data for column O ("tipo")
$ws.Range("O1").Value = "tipo"
$ws.Range("O2").Value = "NÃO"
$ws.Range("O3").Value = "CAPEX"
$ws.Range("O4").Value = "opex"

# Column D grew wider to fit the new content elsewhere in the sheet.
# 64.16666667 compensates for the engine's implicit padding so the
# persisted <col> width lands exactly on 65.
$ws.Range("D1").ColumnWidth = 64.16666666666667

# Scroll/selection bookkeeping, mirroring what Excel records after a user
# edits the new column and leaves the cursor on the last written cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("O5").Select()
